$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "30.245.18"
$ws.Range("E2").Value = "  +1.91%  "
$ws.Range("D3").Value = "2.089.64"
$ws.Range("E3").Value = "  -0.31%  "
$ws.Range("D4").Value = "'1.004"
$ws.Range("E4").Value = "  -0.40%  "
$ws.Range("D5").Value = "'343.29"
$ws.Range("E5").Value = "  -0.07%  "
$ws.Range("D6").Value = "'1.004"
$ws.Range("E6").Value = "  -0.28%  "
$ws.Range("D7").Value = "'0.5209"
$ws.Range("E7").Value = "  +1.16%  "
$ws.Range("D8").Value = "'0.4402"
$ws.Range("D9").Value = "'54.31"
$ws.Range("E9").Value = "  +2.39%  "
$ws.Range("D10").Value = "'0.09316"
$ws.Range("E10").Value = "  +1.25%  "
$ws.Range("E11").Value = "  -0.14%  "
$ws.Range("D12").Value = "'24.73"
$ws.Range("E12").Value = "  -0.70%  "
$ws.Range("D13").Value = "'8.667"
$ws.Range("E13").Value = "  +5.58%  "
$ws.Range("D14").Value = "2.117.34"
$ws.Range("E14").Value = "  +1.13%  "
$ws.Range("D15").Value = "'6.899"
$ws.Range("E15").Value = "  +2.19%  "
$ws.Range("D16").Value = "'101.17"
$ws.Range("E16").Value = "  +1.80%  "
$ws.Range("D17").Value = "'0.00001156"
$ws.Range("E17").Value = "  +0.50%  "
$ws.Range("D18").Value = "'1.005"
$ws.Range("E18").Value = "  -0.35%  "
$ws.Range("D19").Value = "'21.17"
$ws.Range("E19").Value = "  +2.32%  "
$ws.Range("D20").Value = "'0.06690"
$ws.Range("E20").Value = "  +0.87%  "
$ws.Range("D21").Value = "'6.368"
$ws.Range("E21").Value = "  +2.85%  "
$ws.Range("D22").Value = "'1.005"
$ws.Range("E22").Value = "  -0.14%  "
$ws.Range("D23").Value = "30.269.79"
$ws.Range("E23").Value = "  +1.81%  "
$ws.Range("D24").Value = "'12.51"
$ws.Range("E24").Value = "  -0.52%  "
$ws.Range("D25").Value = "'2.289"
$ws.Range("E25").Value = "  -0.81%  "
$ws.Range("D26").Value = "'21.69"
$ws.Range("E26").Value = "  -0.79%  "
$ws.Range("D27").Value = "'162.04"
$ws.Range("E27").Value = "  -0.40%  "
$ws.Range("D28").Value = "'2.516"
$ws.Range("E28").Value = "  -0.34%  "
$ws.Range("D29").Value = "'132.89"
$ws.Range("E29").Value = "  +0.25%  "
$ws.Range("D30").Value = "'1.129"
$ws.Range("E30").Value = "  +0.01%  "
$ws.Range("B31").Value = "ARBITRUM"
$ws.Range("C31").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D31").Value = "'1.668"
$ws.Range("E31").Value = "  +1.47%  "
$ws.Range("B32").Value = "Stellar"
$ws.Range("C32").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D32").Value = "'0.1046"
$ws.Range("E32").Value = "  -0.19%  "
$ws.Range("B33").Value = "Filecoin"
$ws.Range("C33").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D33").Value = "'6.211"
$ws.Range("E33").Value = "  +0.90%  "
$ws.Range("B34").Value = "InternetComputer(DFINITY)"
$ws.Range("C34").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D34").Value = "'6.685"
$ws.Range("E34").Value = "  +10.49%  "
$ws.Range("D35").Value = "'3.861"
$ws.Range("E35").Value = "  -2.55%  "
$ws.Range("E36").Value = "  -1.41%  "
$ws.Range("D37").Value = "'0.02627"
$ws.Range("E37").Value = "  +2.53%  "
$ws.Range("D38").Value = "'0.06747"
$ws.Range("E38").Value = "  +0.54%  "
$ws.Range("B39").Value = "TrustWalletToken"
$ws.Range("C39").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D39").Value = "'1.351"
$ws.Range("E39").Value = "  +4.41%  "
$ws.Range("B40").Value = "TheSandbox"
$ws.Range("C40").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D40").Value = "'0.6975"
$ws.Range("E40").Value = "  +1.67%  "
$ws.Range("D41").Value = "'12.52"
$ws.Range("E41").Value = "  +1.08%  "
$ws.Range("D42").Value = "'0.2211"
$ws.Range("E42").Value = "  -0.51%  "
$ws.Range("D43").Value = "'0.6832"
$ws.Range("E43").Value = "  +2.89%  "
$ws.Range("D44").Value = "'14.37"
$ws.Range("E44").Value = "  +1.60%  "
$ws.Range("D45").Value = "'2.340"
$ws.Range("E45").Value = "  +2.25%  "
$ws.Range("D46").Value = "'1.004"
$ws.Range("E46").Value = "  -0.14%  "
$ws.Range("D47").Value = "'1.359"
$ws.Range("E47").Value = "  +17.17%  "
$ws.Range("D48").Value = "'3.631"
$ws.Range("E48").Value = "  +0.50%  "
$ws.Range("D49").Value = "'0.00000000346"
$ws.Range("E49").Value = "  -0.80%  "
$ws.Range("D50").Value = "'1.217"
$ws.Range("E50").Value = "  +9.48%  "
$ws.Range("D51").Value = "'1.214"
$ws.Range("E51").Value = "  -0.26%  "
